$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats paste-special mode (formats only).
$xlPasteFormats = -4122

# Drop the old last data row (14) first - the refreshed table only has 12
# data rows (2-13) instead of 13. Doing this before any other structural
# change keeps the sheet's used-range/dimension bookkeeping consistent.
$ws.Range("A14:G14").Delete()

# Snapshot the two alternating row styles ("pink" = team1 highlighted,
# "green" = team2 highlighted... actually: pink on the non-favoured team,
# green on the favoured/predicted-winner team) into untouched scratch rows
# far below the table, before we start overwriting rows 2-13 (which would
# otherwise make a live reference to e.g. row 2 or row 5 pick up the new
# values/formats instead of the original ones).
$ws.Range("A2:G2").Copy()
$ws.Range("A200:G200").PasteSpecial($xlPasteFormats)
$ws.Range("A5:G5").Copy()
$ws.Range("A201:G201").PasteSpecial($xlPasteFormats)
$pink = $ws.Range("A200:G200")
$green = $ws.Range("A201:G201")

# --- Row 2 (style pattern 2) ---
$ws.Range("A2").Value = "Dallas Mavericks"
$ws.Range("B2").Value = "Miami Heat"
$ws.Range("C2").Value = -130
$ws.Range("D2").Value = 110
$ws.Range("E2").Value = "Miami Heat"
$ws.Range("F2").Value = 0.288399363478462
$ws.Range("G2").Value = 0.7116006265874305
$pink.Copy()
$ws.Range("A2:G2").PasteSpecial($xlPasteFormats)

# --- Row 3 (style pattern 2) ---
$ws.Range("A3").Value = "Los Angeles Clippers"
$ws.Range("B3").Value = "New Orleans Pelicans"
$ws.Range("C3").Value = 140
$ws.Range("D3").Value = -165
$ws.Range("E3").Value = "New Orleans Pelicans"
$ws.Range("F3").Value = 0.2218646870660309
$ws.Range("G3").Value = 0.7781353029998616
$pink.Copy()
$ws.Range("A3:G3").PasteSpecial($xlPasteFormats)

# --- Row 4 (style pattern 2) ---
$ws.Range("A4").Value = "Charlotte Hornets"
$ws.Range("B4").Value = "Toronto Raptors"
$ws.Range("C4").Value = 600
$ws.Range("D4").Value = -850
$ws.Range("E4").Value = "Toronto Raptors"
$ws.Range("F4").Value = 0.1654025785546778
$ws.Range("G4").Value = 0.8345974115112148
$pink.Copy()
$ws.Range("A4:G4").PasteSpecial($xlPasteFormats)

# --- Row 5 (style pattern 2) ---
$ws.Range("A5").Value = "Chicago Bulls"
$ws.Range("B5").Value = "Memphis Grizzlies"
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = -130
$ws.Range("E5").Value = "Memphis Grizzlies"
$ws.Range("F5").Value = 0.1918798322583874
$ws.Range("G5").Value = 0.8081201578075051
$pink.Copy()
$ws.Range("A5:G5").PasteSpecial($xlPasteFormats)

# --- Row 6 (style pattern 3) ---
$ws.Range("A6").Value = "Minnesota Timberwolves"
$ws.Range("B6").Value = "Portland Trail Blazers"
$ws.Range("C6").Value = -1800
$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = "Minnesota Timberwolves"
$ws.Range("F6").Value = 0.7576219517000373
$ws.Range("G6").Value = 0.2423780482999626
$green.Copy()
$ws.Range("A6:G6").PasteSpecial($xlPasteFormats)

# --- Row 7 (style pattern 3) ---
$ws.Range("A7").Value = "Atlanta Hawks"
$ws.Range("B7").Value = "Dallas Mavericks"
$ws.Range("C7").Value = -150
$ws.Range("D7").Value = 130
$ws.Range("E7").Value = "Atlanta Hawks"
$ws.Range("F7").Value = 0.7607025830436108
$ws.Range("G7").Value = 0.2392974169563893
$green.Copy()
$ws.Range("A7:G7").PasteSpecial($xlPasteFormats)

# --- Row 8 (style pattern 2) ---
$ws.Range("A8").Value = "Detroit Pistons"
$ws.Range("B8").Value = "Orlando Magic"
$ws.Range("C8").Value = 295
$ws.Range("D8").Value = -360
$ws.Range("E8").Value = "Orlando Magic"
$ws.Range("F8").Value = 0.1751287282125945
$ws.Range("G8").Value = 0.8248712618532981
$pink.Copy()
$ws.Range("A8:G8").PasteSpecial($xlPasteFormats)

# --- Row 9 (style pattern 3) ---
$ws.Range("A9").Value = "Sacramento Kings"
$ws.Range("B9").Value = "San Antonio Spurs"
$ws.Range("C9").Value = -1400
$ws.Range("D9").Value = 850
$ws.Range("E9").Value = "Sacramento Kings"
$ws.Range("F9").Value = 0.7769561269067223
$ws.Range("G9").Value = 0.2230438730932777
$green.Copy()
$ws.Range("A9:G9").PasteSpecial($xlPasteFormats)

# --- Row 10 (style pattern 2) ---
$ws.Range("A10").Value = "Houston Rockets"
$ws.Range("B10").Value = "Los Angeles Lakers"
$ws.Range("C10").Value = 470
$ws.Range("D10").Value = -625
$ws.Range("E10").Value = "Los Angeles Lakers"
$ws.Range("F10").Value = 0.1384317064644886
$ws.Range("G10").Value = 0.8615682836014039
$pink.Copy()
$ws.Range("A10:G10").PasteSpecial($xlPasteFormats)

# --- Row 11 (style pattern 2) ---
$ws.Range("A11").Value = "Oklahoma City Thunder"
$ws.Range("B11").Value = "Phoenix Suns"
$ws.Range("C11").Value = 170
$ws.Range("D11").Value = -200
$ws.Range("E11").Value = "Phoenix Suns"
$ws.Range("F11").Value = 0.3359166330691408
$ws.Range("G11").Value = 0.6640833669308592
$pink.Copy()
$ws.Range("A11:G11").PasteSpecial($xlPasteFormats)

# --- Row 12 (style pattern 3) ---
$ws.Range("A12").Value = "Cleveland Cavaliers"
$ws.Range("B12").Value = "Indiana Pacers"
$ws.Range("C12").Value = -750
$ws.Range("D12").Value = 550
$ws.Range("E12").Value = "Cleveland Cavaliers"
$ws.Range("F12").Value = 0.6431379835304963
$ws.Range("G12").Value = 0.3568620065353962
$green.Copy()
$ws.Range("A12:G12").PasteSpecial($xlPasteFormats)

# --- Row 13 (style pattern 3) ---
$ws.Range("A13").Value = "Milwaukee Bucks"
$ws.Range("B13").Value = "Philadelphia 76ers"
$ws.Range("C13").Value = -190
$ws.Range("D13").Value = 160
$ws.Range("E13").Value = "Milwaukee Bucks"
$ws.Range("F13").Value = 0.7653531779564172
$ws.Range("G13").Value = 0.2346468220435829
$green.Copy()
$ws.Range("A13:G13").PasteSpecial($xlPasteFormats)


# Clean up the scratch rows used as style templates.
$ws.Range("A200:G201").Delete()
